$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new column header: AC1 = c_thouy ---
$ws.Cells.Item(1,29).Value = "c_thouy"

# --- corrected sample sizes for existing rows 134 / 135 ---
$ws.Cells.Item(134,9).Value = 1830
$ws.Cells.Item(135,9).Value = 745

# --- row 136 ---
$ws.Cells.Item(136,1).Value = 42
$ws.Cells.Item(136,2).Value = 2021
$ws.Cells.Item(136,3).Value = 15
$ws.Cells.Item(136,4).Value = 12
$ws.Cells.Item(136,5).Value = 7
$ws.Cells.Item(136,6).Value = "elabe"
$ws.Cells.Item(136,7).Value = "online"
$ws.Cells.Item(136,8).Value = "partially"
$ws.Cells.Item(136,9).Value = 934
$ws.Cells.Item(136,10).Value = 2
$ws.Cells.Item(136,11).Value = 1
$ws.Cells.Item(136,12).Value = 8
$ws.Cells.Item(136,13).Value = 1
$ws.Cells.Item(136,14).Value = 2
$ws.Cells.Item(136,15).Value = 7
$ws.Cells.Item(136,16).Value = 3
$ws.Cells.Item(136,17).Value = 23
$ws.Cells.Item(136,18).Value = 20
$ws.Cells.Item(136,21).Value = 2
$ws.Cells.Item(136,22).Value = 2
$ws.Cells.Item(136,23).Value = 15
$ws.Cells.Item(136,24).Value = 13
$ws.Cells.Item(136,25).Value = "T_1"
$ws.Cells.Item(136,25).Font.Color = 0
$ws.Cells.Item(136,27).Value = "T_1"
$ws.Cells.Item(136,29).Value = "T_1"

# --- row 137 ---
$ws.Cells.Item(137,1).Value = 43
$ws.Cells.Item(137,2).Value = 2021
$ws.Cells.Item(137,3).Value = 15
$ws.Cells.Item(137,4).Value = 12
$ws.Cells.Item(137,5).Value = 7
$ws.Cells.Item(137,6).Value = "ipsos"
$ws.Cells.Item(137,7).Value = "online"
$ws.Cells.Item(137,8).Value = "excluded"
$ws.Cells.Item(137,9).Value = 747
$ws.Cells.Item(137,10).Value = 1.5
$ws.Cells.Item(137,11).Value = 0.5
$ws.Cells.Item(137,12).Value = 8
$ws.Cells.Item(137,13).Value = 2
$ws.Cells.Item(137,14).Value = 2.5
$ws.Cells.Item(137,15).Value = 7
$ws.Cells.Item(137,16).Value = 5
$ws.Cells.Item(137,17).Value = 25
$ws.Cells.Item(137,18).Value = 16
$ws.Cells.Item(137,21).Value = 1
$ws.Cells.Item(137,22).Value = 1.5
$ws.Cells.Item(137,23).Value = 16
$ws.Cells.Item(137,24).Value = 14

# --- row 138 ---
$ws.Cells.Item(138,1).Value = 44
$ws.Cells.Item(138,2).Value = 2021
$ws.Cells.Item(138,3).Value = 15
$ws.Cells.Item(138,4).Value = 12
$ws.Cells.Item(138,5).Value = 7
$ws.Cells.Item(138,6).Value = "bva"
$ws.Cells.Item(138,7).Value = "online"
$ws.Cells.Item(138,8).Value = "excluded"
$ws.Cells.Item(138,9).Value = 894
$ws.Cells.Item(138,10).Value = 1.5
$ws.Cells.Item(138,11).Value = 0.5
$ws.Cells.Item(138,12).Value = 9
$ws.Cells.Item(138,13).Value = 2.5
$ws.Cells.Item(138,14).Value = 1
$ws.Cells.Item(138,15).Value = 7
$ws.Cells.Item(138,16).Value = 5
$ws.Cells.Item(138,17).Value = 24
$ws.Cells.Item(138,18).Value = 17
$ws.Cells.Item(138,21).Value = 1
$ws.Cells.Item(138,22).Value = 2.5
$ws.Cells.Item(138,23).Value = 16
$ws.Cells.Item(138,24).Value = 13
$ws.Cells.Item(138,25).Value = "T_0.5"
$ws.Cells.Item(138,25).Font.Color = 0

# --- row 139 ---
$ws.Cells.Item(139,1).Value = 45
$ws.Cells.Item(139,2).Value = 2021
$ws.Cells.Item(139,3).Value = 15
$ws.Cells.Item(139,4).Value = 12
$ws.Cells.Item(139,5).Value = 8
$ws.Cells.Item(139,6).Value = "odoxa"
$ws.Cells.Item(139,7).Value = "online"
$ws.Cells.Item(139,8).Value = "partially"
$ws.Cells.Item(139,9).Value = 1391
$ws.Cells.Item(139,10).Value = 1.5
$ws.Cells.Item(139,11).Value = 1
$ws.Cells.Item(139,12).Value = 10
$ws.Cells.Item(139,13).Value = 2
$ws.Cells.Item(139,14).Value = 1
$ws.Cells.Item(139,15).Value = 6
$ws.Cells.Item(139,16).Value = 3
$ws.Cells.Item(139,17).Value = 24
$ws.Cells.Item(139,18).Value = 19
$ws.Cells.Item(139,21).Value = 1
$ws.Cells.Item(139,22).Value = 2.5
$ws.Cells.Item(139,23).Value = 17
$ws.Cells.Item(139,24).Value = 12

# --- row 140 ---
$ws.Cells.Item(140,1).Value = 46
$ws.Cells.Item(140,2).Value = 2021
$ws.Cells.Item(140,3).Value = 15
$ws.Cells.Item(140,4).Value = 12
$ws.Cells.Item(140,5).Value = 6
$ws.Cells.Item(140,6).Value = "cluster17"
$ws.Cells.Item(140,7).Value = "online"
$ws.Cells.Item(140,8).Value = "partially"
$ws.Cells.Item(140,9).Value = 1487
$ws.Cells.Item(140,10).Value = 1.5
$ws.Cells.Item(140,11).Value = 0.5
$ws.Cells.Item(140,12).Value = 13
$ws.Cells.Item(140,13).Value = 2
$ws.Cells.Item(140,14).Value = 1
$ws.Cells.Item(140,15).Value = 5
$ws.Cells.Item(140,16).Value = 3
$ws.Cells.Item(140,17).Value = 23
$ws.Cells.Item(140,18).Value = 16
$ws.Cells.Item(140,21).Value = 1
$ws.Cells.Item(140,22).Value = 2
$ws.Cells.Item(140,23).Value = 16
$ws.Cells.Item(140,24).Value = 15
$ws.Cells.Item(140,25).Value = 0.5
$ws.Cells.Item(140,27).Value = 1

# --- final selection, matching the saved view state ---
[void]$ws.Range("AA140").Select()

